$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.261.78'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '2.719.90'
$ws.Range("E3").Value = '  +2.30%  '
$ws.Range("E4").Value = '  -0.05%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '611.75'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +0.90%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.25'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  +1.36%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.57%  '
$ws.Range("E9").Value = '  +6.82%  '
$ws.Range("E10").Value = '  +3.69%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.406'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("E12").Value = '  +1.55%  '
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '30.47'
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  +4.16%  '
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000210'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  +16.00%  '
$ws.Range("D15").Value = '3.208.03'
$ws.Range("E15").Value = '  +2.31%  '
$ws.Range("D16").Value = '66.074.24'
$ws.Range("E16").Value = '  +1.21%  '
$ws.Range("D17").Value = '2.724.67'
$ws.Range("E17").Value = '  -5.21%  '
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.86'
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = '  +1.60%  '
$ws.Range("E19").Value = '  +1.50%  '
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '363.03'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  +2.08%  '
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.68'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  +4.67%  '
$ws.Range("E22").Value = '  -0.03%  '
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.63'
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  +3.25%  '
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.84'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  +2.69%  '
$ws.Range("E25").Value = '  +12.54%  '
$ws.Range("E26").Value = '  -2.99%  '
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.74'
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = '  +3.67%  '
$ws.Range("E28").Value = '  +4.44%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("E30").Value = '  +4.13%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '540.73'
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("E33").Value = '  +0.35%  '
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.74'
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  +4.05%  '
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.51'
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = '  -5.14%  '
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.438'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  +1.73%  '
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '21.04'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  +3.78%  '
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.67'
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("E39").Value = '  -2.50%  '
$ws.Range("E40").Value = '  -0.11%  '
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '171.62'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  +1.98%  '
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.21'
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = '  +1.80%  '
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0621'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  +1.15%  '
$ws.Range("E46").Value = '  +2.91%  '
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.76'
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = '  +1.53%  '
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0268'
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = '  +5.42%  '
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.39'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  +8.18%  '
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0991'
$ws.Range("D51").Style = $style
$ws.Range("E51").Value = '  +0.56%  '
